$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This handback-status report gains one new row for the file
# "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" across all three sheets
# (Overview, zh-cn, de-de). The new row is inserted right before the
# existing "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" row, pushing that one
# down by one row on every sheet and growing each table by one row.
#
# NB: custom functions here only bind POSITIONAL parameters reliably, so
# everything below is passed positionally.
# ---------------------------------------------------------------------------

function Add-Row3AndFixHyperlinks(
    $SheetName,
    $TableRangeA1,
    $NewRowValues,
    $LinkCols,
    $NewFileDisplay,
    $NewFileUrl,
    $OldFileDisplay,
    $OldFileUrl,
    $FirstFileDisplay,
    $FirstFileUrl
) {
    $ws = $wb.Worksheets.Item($SheetName)
    $lo = $ws.ListObjects.Item(1)

    # Push the "e697e9da" row (currently row 3) down to row 4, carrying
    # styles/number-formats with it.
    $ws.Rows(3).Insert()

    # Grow the table / autofilter to include the new row.
    $lo.Resize($ws.Range($TableRangeA1))

    # Populate the freshly-inserted row 3 with the new record's data.
    foreach ($col in $NewRowValues.Keys) {
        $ws.Range($col + "3").Value = $NewRowValues[$col]
    }

    # The hyperlink ranges don't auto-shift with the row insert, and the
    # engine's Range.Hyperlinks.Delete() drops every hyperlink on the sheet,
    # so rebuild the whole collection in top-to-bottom order.
    $ws.Range("A1").Hyperlinks.Delete()

    foreach ($col in $LinkCols) {
        $ws.Hyperlinks.Add($ws.Range($col + "2"), $FirstFileUrl, "", "", $FirstFileDisplay) | Out-Null
    }
    foreach ($col in $LinkCols) {
        $ws.Hyperlinks.Add($ws.Range($col + "3"), $NewFileUrl, "", "", $NewFileDisplay) | Out-Null
    }
    foreach ($col in $LinkCols) {
        $ws.Hyperlinks.Add($ws.Range($col + "4"), $OldFileUrl, "", "", $OldFileDisplay) | Out-Null
    }
}

# --------------------------- Overview sheet -------------------------------
$overviewValues = @{
    "A" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "B" = "e2e\d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "C" = ".md"
    "E" = "Handed back: in sync with en-US"
    "F" = "Handed back: in sync with en-US"
    "G" = "2016-08-20 12:48:01"
}
Add-Row3AndFixHyperlinks `
    "Overview" `
    "A1:G4" `
    $overviewValues `
    @("B") `
    "e2e\d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2093b342de8e20a28f3cb96610e78449bd46cfb3/e2e/d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "e2e\e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2093b342de8e20a28f3cb96610e78449bd46cfb3/e2e/e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "e2e\df6a600a-6246-4366-bdc7-3466dd5b0682.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a980daef40b2f38ea1a35de72c0329219b2c258b/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md"

# ------------------------------ zh-cn sheet --------------------------------
$zhcnValues = @{
    "A" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "B" = ".md"
    "C" = "Handed back: in sync with en-US"
    "D" = "e2e"
    "E" = "ht"
    "F" = "True"
    "G" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.8dea39b6f0615507ba1cb6d8ac9006756ec339d5.zh-cn.xlf"
    "H" = "2016-08-20 12:47:56"
    "I" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "J" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.8dea39b6f0615507ba1cb6d8ac9006756ec339d5.zh-cn.xlf"
    "K" = "2016-08-20 12:48:24"
    "L" = ""
    "M" = "True"
    "N" = ""
    "O" = "False"
    "P" = ""
}
Add-Row3AndFixHyperlinks `
    "zh-cn" `
    "A1:P4" `
    $zhcnValues `
    @("A", "I") `
    "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9b348f2a50bd51b632f403a85bf8315bb6972a1/e2e/d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4e0f1a73aa74fe0c8cde6f92ed400a01c73371c4/e2e/e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "df6a600a-6246-4366-bdc7-3466dd5b0682.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f280221a50bd51b632f403a85bf8315bb6972a12/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md"

# ------------------------------ de-de sheet --------------------------------
$dedeValues = @{
    "A" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "B" = ".md"
    "C" = "Handed back: in sync with en-US"
    "D" = "e2e"
    "E" = "ht"
    "F" = "True"
    "G" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.8dea39b6f0615507ba1cb6d8ac9006756ec339d5.de-de.xlf"
    "H" = "2016-08-20 12:48:01"
    "I" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md"
    "J" = "d9b348f2-6622-4b99-9ea5-727a219a8dd9.8dea39b6f0615507ba1cb6d8ac9006756ec339d5.de-de.xlf"
    "K" = "2016-08-20 12:48:31"
    "L" = ""
    "M" = "True"
    "N" = ""
    "O" = "False"
    "P" = ""
}
Add-Row3AndFixHyperlinks `
    "de-de" `
    "A1:P4" `
    $dedeValues `
    @("A", "I") `
    "d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d9b348f2e9c2134c1e8e27d3a647feb27c36adf6/e2e/d9b348f2-6622-4b99-9ea5-727a219a8dd9.md" `
    "e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a054f43409e5be2e985bd3faddb411ed2b4220dc/e2e/e697e9da-d9a5-4a6b-8522-4a5a26a8a6de.md" `
    "df6a600a-6246-4366-bdc7-3466dd5b0682.md" `
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bf3632aed9c2134c1e8e27d3a647feb27c36adf6/e2e/df6a600a-6246-4366-bdc7-3466dd5b0682.md"
